# Scheduled runner: refresh market-derived columns (H:N) across all Leve sheets
# with latest pricing data. Pure value overwrite, no formulas/structure changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1988.7778
$ws.Range("J41").Value = 2112.375
$ws.Range("L41").Value = 2112.375
$ws.Range("N41").Value = -2992.375
$ws.Range("H74").Value = 4099.6
$ws.Range("I74").Value = 2749.5
$ws.Range("K74").Value = 2749.5
$ws.Range("M74").Value = -1813.5
$ws.Range("H77").Value = 4099.6
$ws.Range("I77").Value = 2749.5
$ws.Range("K77").Value = 13747.5
$ws.Range("M77").Value = -9067.5
$ws.Range("H98").Value = 1428.8064
$ws.Range("I98").Value = 1037.75
$ws.Range("J98").Value = 2769.5715
$ws.Range("K98").Value = 1037.75
$ws.Range("L98").Value = 2769.5715
$ws.Range("M98").Value = 460.25
$ws.Range("N98").Value = -5765.5715
$ws.Range("H100").Value = 2514
$ws.Range("I100").Value = 1419.2
$ws.Range("J100").Value = 3882.5
$ws.Range("K100").Value = 1419.2
$ws.Range("L100").Value = 3882.5
$ws.Range("M100").Value = -878.2
$ws.Range("N100").Value = -4964.5
$ws.Range("H116").Value = 16273.875
$ws.Range("J116").Value = 4313
$ws.Range("L116").Value = 4313
$ws.Range("N116").Value = -11197
$ws.Range("H122").Value = 1428.8064
$ws.Range("I122").Value = 1037.75
$ws.Range("J122").Value = 2769.5715
$ws.Range("K122").Value = 3113.25
$ws.Range("L122").Value = 8308.7145
$ws.Range("M122").Value = -663.25
$ws.Range("N122").Value = -13208.7145
$ws.Range("H132").Value = 816.6875
$ws.Range("I132").Value = 715.0714
$ws.Range("J132").Value = 1528
$ws.Range("K132").Value = 2145.2142
$ws.Range("L132").Value = 4584
$ws.Range("M132").Value = 384.7857999999997
$ws.Range("N132").Value = -9644
$ws.Range("H138").Value = 1813.1552
$ws.Range("I138").Value = 1296.64
$ws.Range("J138").Value = 2204.4546
$ws.Range("K138").Value = 3889.92
$ws.Range("L138").Value = 6613.3638
$ws.Range("M138").Value = 1250.08
$ws.Range("N138").Value = -16893.3638
$ws.Range("H141").Value = 5602378
$ws.Range("I141").Value = 7001723
$ws.Range("K141").Value = 21005169
$ws.Range("M141").Value = -20999989

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5676.673
$ws.Range("I32").Value = 3585.8696
$ws.Range("K32").Value = 3585.8696
$ws.Range("M32").Value = -3298.8696
$ws.Range("H61").Value = 6590.773
$ws.Range("I61").Value = 7907.769
$ws.Range("J61").Value = 4688.4443
$ws.Range("K61").Value = 7907.769
$ws.Range("L61").Value = 4688.4443
$ws.Range("M61").Value = -7695.769
$ws.Range("N61").Value = -5112.4443
$ws.Range("H74").Value = 3549.5
$ws.Range("I74").Value = 1500
$ws.Range("J74").Value = 3959.4
$ws.Range("K74").Value = 1500
$ws.Range("L74").Value = 3959.4
$ws.Range("M74").Value = -626
$ws.Range("N74").Value = -5707.4
$ws.Range("H77").Value = 3549.5
$ws.Range("I77").Value = 1500
$ws.Range("J77").Value = 3959.4
$ws.Range("K77").Value = 7500
$ws.Range("L77").Value = 19797
$ws.Range("M77").Value = -3132
$ws.Range("N77").Value = -28533
$ws.Range("H110").Value = 1735.8572
$ws.Range("I110").Value = 300.33334
$ws.Range("K110").Value = 300.33334
$ws.Range("M110").Value = 1744.66666
$ws.Range("H132").Value = 1620.4242
$ws.Range("I132").Value = 1278.52
$ws.Range("K132").Value = 3835.56
$ws.Range("M132").Value = -1305.56
$ws.Range("H135").Value = 34724.75
$ws.Range("J135").Value = 34724.75
$ws.Range("L135").Value = 34724.75
$ws.Range("N135").Value = -44864.75
$ws.Range("H136").Value = 6590.773
$ws.Range("I136").Value = 7907.769
$ws.Range("J136").Value = 4688.4443
$ws.Range("K136").Value = 23723.307
$ws.Range("L136").Value = 14065.3329
$ws.Range("M136").Value = -21173.307
$ws.Range("N136").Value = -19165.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5000
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5494
$ws.Range("H105").Value = 2368.2
$ws.Range("J105").Value = 3500
$ws.Range("L105").Value = 3500
$ws.Range("N105").Value = -6994
$ws.Range("H107").Value = 971.9286
$ws.Range("I107").Value = 656.3333
$ws.Range("J107").Value = 1208.625
$ws.Range("K107").Value = 656.3333
$ws.Range("L107").Value = 1208.625
$ws.Range("M107").Value = 1263.6667
$ws.Range("N107").Value = -5048.625
$ws.Range("H134").Value = 5915.88
$ws.Range("I134").Value = 6598.7
$ws.Range("K134").Value = 19796.1
$ws.Range("M134").Value = -17261.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 500
$ws.Range("I35").Value = 500
$ws.Range("K35").Value = 500
$ws.Range("M35").Value = -206
$ws.Range("H38").Value = 1699.5
$ws.Range("I38").Value = 1899
$ws.Range("J38").Value = 1500
$ws.Range("K38").Value = 1899
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = -1522
$ws.Range("N38").Value = -2254
$ws.Range("H46").Value = 1699.5
$ws.Range("I46").Value = 1899
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 1899
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -1688
$ws.Range("N46").Value = -1922

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 260
$ws.Range("I4").Value = 175
$ws.Range("K4").Value = 525
$ws.Range("M4").Value = -413
$ws.Range("H5").Value = 400.78946
$ws.Range("I5").Value = 338.5625
$ws.Range("J5").Value = 732.6667
$ws.Range("K5").Value = 1015.6875
$ws.Range("L5").Value = 2198.0001
$ws.Range("M5").Value = -903.6875
$ws.Range("N5").Value = -2422.0001
$ws.Range("H17").Value = 20251.5
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H126").Value = 5996.6665
$ws.Range("I126").Value = 5990
$ws.Range("K126").Value = 17970
$ws.Range("M126").Value = -13030
$ws.Range("H129").Value = 121366
$ws.Range("J129").Value = 241822.67
$ws.Range("L129").Value = 725468.01
$ws.Range("N129").Value = -735468.01
$ws.Range("H130").Value = 2421.1667
$ws.Range("J130").Value = 3499.6667
$ws.Range("L130").Value = 10499.0001
$ws.Range("N130").Value = -20539.0001
$ws.Range("H131").Value = 11267.513
$ws.Range("J131").Value = 11551.365
$ws.Range("L131").Value = 34654.095
$ws.Range("N131").Value = -44734.095
$ws.Range("H135").Value = 400.78946
$ws.Range("I135").Value = 338.5625
$ws.Range("J135").Value = 732.6667
$ws.Range("K135").Value = 3047.0625
$ws.Range("L135").Value = 6594.0003
$ws.Range("M135").Value = -512.0625
$ws.Range("N135").Value = -11664.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 125.833336
$ws.Range("I2").Value = 80.75
$ws.Range("J2").Value = 216
$ws.Range("K2").Value = 80.75
$ws.Range("L2").Value = 216
$ws.Range("M2").Value = 32.25
$ws.Range("N2").Value = -442
$ws.Range("H10").Value = 507499
$ws.Range("I10").Value = 673332.3
$ws.Range("K10").Value = 673332.3
$ws.Range("M10").Value = -673163.3
$ws.Range("H14").Value = 2026251.1
$ws.Range("I14").Value = 3041599.8
$ws.Range("J14").Value = 334003.34
$ws.Range("K14").Value = 3041599.8
$ws.Range("L14").Value = 334003.34
$ws.Range("M14").Value = -3041431.8
$ws.Range("N14").Value = -334339.34
$ws.Range("H70").Value = 3949.75
$ws.Range("I70").Value = 3800
$ws.Range("K70").Value = 3800
$ws.Range("M70").Value = -3530
$ws.Range("H73").Value = 3949.75
$ws.Range("I73").Value = 3800
$ws.Range("K73").Value = 3800
$ws.Range("M73").Value = -2864
$ws.Range("H107").Value = 96.666664
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 1396
$ws.Range("I113").Value = 1135.5
$ws.Range("K113").Value = 1135.5
$ws.Range("M113").Value = 1034.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2348.5789
$ws.Range("I7").Value = 2242.5293
$ws.Range("K7").Value = 2242.5293
$ws.Range("M7").Value = -2130.5293
$ws.Range("H40").Value = 11643.75
$ws.Range("I40").Value = 12218.728
$ws.Range("J40").Value = 10378.8
$ws.Range("K40").Value = 12218.728
$ws.Range("L40").Value = 10378.8
$ws.Range("M40").Value = -12082.728
$ws.Range("N40").Value = -10650.8
$ws.Range("H46").Value = 1613.0588
$ws.Range("I46").Value = 985.5454999999999
$ws.Range("K46").Value = 985.5454999999999
$ws.Range("M46").Value = -797.5454999999999
$ws.Range("H122").Value = 7564.2144
$ws.Range("I122").Value = 7099.3335
$ws.Range("J122").Value = 8401
$ws.Range("K122").Value = 21298.0005
$ws.Range("L122").Value = 25203
$ws.Range("M122").Value = -18848.0005
$ws.Range("N122").Value = -30103
$ws.Range("H126").Value = 2348.5789
$ws.Range("I126").Value = 2242.5293
$ws.Range("K126").Value = 6727.5879
$ws.Range("M126").Value = -4257.5879
$ws.Range("H132").Value = 1845.2222
$ws.Range("I132").Value = 1265.1
$ws.Range("K132").Value = 3795.3
$ws.Range("M132").Value = -1265.3
$ws.Range("H136").Value = 4733.222
$ws.Range("I136").Value = 3265.6667
$ws.Range("K136").Value = 9797.000100000001
$ws.Range("M136").Value = -7247.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1052.4062
$ws.Range("I132").Value = 826.69385
$ws.Range("K132").Value = 2480.08155
$ws.Range("M132").Value = 49.91845000000012
$ws.Range("H136").Value = 22223636
$ws.Range("J136").Value = 1784.8
$ws.Range("L136").Value = 5354.4
$ws.Range("N136").Value = -10454.4
